# Navigation and switchable layers
# Reorders the "layers" sheet so prevnav/nextnav become the first
# interactive layers (levels 4-5), followed by menu, leftmenu, maintitle,
# titles and topmenu. Also updates menu's alignment/visibility flag and
# leftmenu's source path, and moves the sheet's active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("layers")

# --- Row 6 (level 4): was "menu", now "prevnav" ---
$ws.Range("B6").Value = "prevnav"
$ws.Range("C6").Value = "{templatePath}prevnav"
$ws.Range("D6").Value = 1
$ws.Range("E6").ClearContents()
$ws.Range("F6").Value = "left"
$ws.Range("G6").Value = "l70"
$ws.Range("H6").Value = 30
$ws.Range("I6").Value = 100
$ws.Range("J6").ClearContents()

# --- Row 7 (level 5): was "leftmenu", now "nextnav" ---
$ws.Range("B7").Value = "nextnav"
$ws.Range("C7").Value = "{templatePath}nextnav"
$ws.Range("D7").Value = 1
$ws.Range("F7").Value = "right"
$ws.Range("G7").ClearContents()
$ws.Range("H7").Value = 30
$ws.Range("I7").Value = 100

# --- Row 8 (level 6): was "maintitle", now "menu" ---
$ws.Range("B8").Value = "menu"
$ws.Range("C8").Value = "{templatePath}menu"
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = "top,left"
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("J8").ClearContents()

# --- Row 9 (level 7): was "titles", now "leftmenu" ---
$ws.Range("B9").Value = "leftmenu"
$ws.Range("C9").Value = "{templatePath}leftmenu"
$ws.Range("D9").Value = 1
$ws.Range("F9").Value = "top,left"
$ws.Range("H9").Value = 70
$ws.Range("I9").ClearContents()

# --- Row 10 (level 8): was "topmenu", now "maintitle" ---
$ws.Range("B10").Value = "maintitle"
$ws.Range("C10").ClearContents()
$ws.Range("D10").Value = 1
$ws.Range("F10").Value = "top,left"
$ws.Range("I10").Value = 70

# --- Row 11 (level 9): was "prevnav", now "titles" ---
$ws.Range("B11").Value = "titles"
$ws.Range("C11").Value = "{templatePath}titles"
$ws.Range("D11").Value = 1
$ws.Range("F11").Value = "top,left"
$ws.Range("G11").ClearContents()
$ws.Range("H11").ClearContents()
$ws.Range("I11").Value = 70

# --- Row 12 (level 10): was "nextnav", now "topmenu" ---
$ws.Range("B12").Value = "topmenu"
$ws.Range("C12").Value = "{templatePath}topmenu"
$ws.Range("D12").Value = 1
$ws.Range("F12").Value = "top,left"
$ws.Range("H12").ClearContents()
$ws.Range("I12").Value = 70

# --- Move the active selection on the sheet ---
$ws.Activate()
$ws.Range("B14").Select()
